$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 207,8
$data[0,0] = 'Estados Unidos'
$data[0,1] = 245442
$data[0,2] = 565
$data[0,3] = 10411
$data[0,4] = 228933
$data[0,5] = 5421
$data[0,6] = 28
$data[0,7] = 6098
$data[1,0] = 'España'
$data[1,1] = 117710
$data[1,2] = 5645
$data[1,3] = 30513
$data[1,4] = 76262
$data[1,5] = 6416
$data[1,6] = 587
$data[1,7] = 10935
$data[2,0] = 'Italia'
$data[2,1] = 115242
$data[2,2] = 0
$data[2,3] = 18278
$data[2,4] = 83049
$data[2,5] = 4053
$data[2,6] = 0
$data[2,7] = 13915
$data[3,0] = 'Alemania'
$data[3,1] = 85063
$data[3,2] = 269
$data[3,3] = 22440
$data[3,4] = 61512
$data[3,5] = 3936
$data[3,6] = 4
$data[3,7] = 1111
$data[4,0] = 'China'
$data[4,1] = 81620
$data[4,2] = 31
$data[4,3] = 76571
$data[4,4] = 1727
$data[4,5] = 379
$data[4,6] = 4
$data[4,7] = 3322
$data[5,0] = 'Francia'
$data[5,1] = 59105
$data[5,2] = 0
$data[5,3] = 12428
$data[5,4] = 41290
$data[5,5] = 6399
$data[5,6] = 0
$data[5,7] = 5387
$data[6,0] = 'Iran'
$data[6,1] = 53183
$data[6,2] = 2715
$data[6,3] = 17935
$data[6,4] = 31954
$data[6,5] = 4035
$data[6,6] = 134
$data[6,7] = 3294
$data[7,0] = 'Reino Unido'
$data[7,1] = 33718
$data[7,2] = 0
$data[7,3] = 135
$data[7,4] = 30662
$data[7,5] = 163
$data[7,6] = 0
$data[7,7] = 2921
$data[8,0] = 'Suiza'
$data[8,1] = 19303
$data[8,2] = 476
$data[8,3] = 4846
$data[8,4] = 13884
$data[8,5] = 348
$data[8,6] = 37
$data[8,7] = 573
$data[9,0] = 'Turquia'
$data[9,1] = 18135
$data[9,2] = 0
$data[9,3] = 415
$data[9,4] = 17364
$data[9,5] = 1101
$data[9,6] = 0
$data[9,7] = 356
$data[10,0] = 'Belgica'
$data[10,1] = 16770
$data[10,2] = 1422
$data[10,3] = 2872
$data[10,4] = 12755
$data[10,5] = 1205
$data[10,6] = 132
$data[10,7] = 1143
$data[11,0] = 'Paises Bajos'
$data[11,1] = 14697
$data[11,2] = 0
$data[11,3] = 250
$data[11,4] = 13108
$data[11,5] = 1053
$data[11,6] = 0
$data[11,7] = 1339
$data[12,0] = 'Austria'
$data[12,1] = 11301
$data[12,2] = 172
$data[12,3] = 2022
$data[12,4] = 9111
$data[12,5] = 245
$data[12,6] = 10
$data[12,7] = 168
$data[13,0] = 'Canada'
$data[13,1] = 11283
$data[13,2] = 0
$data[13,3] = 1979
$data[13,4] = 9131
$data[13,5] = 120
$data[13,6] = 0
$data[13,7] = 173
$data[14,0] = 'Corea del Sur'
$data[14,1] = 10062
$data[14,2] = 86
$data[14,3] = 6021
$data[14,4] = 3867
$data[14,5] = 55
$data[14,6] = 5
$data[14,7] = 174
$data[15,0] = 'Portugal'
$data[15,1] = 9034
$data[15,2] = 0
$data[15,3] = 68
$data[15,4] = 8757
$data[15,5] = 230
$data[15,6] = 0
$data[15,7] = 209
$data[16,0] = 'Brasil'
$data[16,1] = 8066
$data[16,2] = 22
$data[16,3] = 127
$data[16,4] = 7612
$data[16,5] = 296
$data[16,6] = 3
$data[16,7] = 327
$data[17,0] = 'Israel'
$data[17,1] = 7030
$data[17,2] = 173
$data[17,3] = 338
$data[17,4] = 6655
$data[17,5] = 115
$data[17,6] = 1
$data[17,7] = 37
$data[18,0] = 'Suecia'
$data[18,1] = 5568
$data[18,2] = 0
$data[18,3] = 103
$data[18,4] = 5157
$data[18,5] = 429
$data[18,6] = 0
$data[18,7] = 308
$data[19,0] = 'Australia'
$data[19,1] = 5350
$data[19,2] = 36
$data[19,3] = 585
$data[19,4] = 4737
$data[19,5] = 50
$data[19,6] = 3
$data[19,7] = 28
$data[20,0] = 'Noruega'
$data[20,1] = 5255
$data[20,2] = 108
$data[20,3] = 32
$data[20,4] = 5169
$data[20,5] = 96
$data[20,6] = 4
$data[20,7] = 54
$data[21,0] = 'Rusia'
$data[21,1] = 4149
$data[21,2] = 601
$data[21,3] = 281
$data[21,4] = 3834
$data[21,5] = 8
$data[21,6] = 4
$data[21,7] = 34
$data[22,0] = 'Chequia'
$data[22,1] = 3869
$data[22,2] = 11
$data[22,3] = 71
$data[22,4] = 3752
$data[22,5] = 77
$data[22,6] = 2
$data[22,7] = 46
$data[23,0] = 'Irlanda'
$data[23,1] = 3849
$data[23,2] = 0
$data[23,3] = 5
$data[23,4] = 3746
$data[23,5] = 109
$data[23,6] = 0
$data[23,7] = 98
$data[24,0] = 'Dinamarca'
$data[24,1] = 3672
$data[24,2] = 286
$data[24,3] = 1193
$data[24,4] = 2340
$data[24,5] = 153
$data[24,6] = 16
$data[24,7] = 139
$data[25,0] = 'Chile'
$data[25,1] = 3404
$data[25,2] = 0
$data[25,3] = 335
$data[25,4] = 3051
$data[25,5] = 31
$data[25,6] = 0
$data[25,7] = 18
$data[26,0] = 'Malasia'
$data[26,1] = 3333
$data[26,2] = 217
$data[26,3] = 827
$data[26,4] = 2453
$data[26,5] = 108
$data[26,6] = 3
$data[26,7] = 53
$data[27,0] = 'Rumania'
$data[27,1] = 3183
$data[27,2] = 445
$data[27,3] = 283
$data[27,4] = 2782
$data[27,5] = 83
$data[27,6] = 3
$data[27,7] = 118
$data[28,0] = 'Ecuador'
$data[28,1] = 3163
$data[28,2] = 0
$data[28,3] = 65
$data[28,4] = 2978
$data[28,5] = 100
$data[28,6] = 0
$data[28,7] = 120
$data[29,0] = 'Polonia'
$data[29,1] = 3149
$data[29,2] = 203
$data[29,3] = 56
$data[29,4] = 3034
$data[29,5] = 50
$data[29,6] = 2
$data[29,7] = 59
$data[30,0] = 'Filipinas'
$data[30,1] = 3018
$data[30,2] = 385
$data[30,3] = 52
$data[30,4] = 2830
$data[30,5] = 1
$data[30,6] = 29
$data[30,7] = 136
$data[31,0] = 'Japon'
$data[31,1] = 2617
$data[31,2] = 0
$data[31,3] = 514
$data[31,4] = 2040
$data[31,5] = 60
$data[31,6] = 0
$data[31,7] = 63
$data[32,0] = 'India'
$data[32,1] = 2567
$data[32,2] = 24
$data[32,3] = 192
$data[32,4] = 2303
$data[32,5] = 0
$data[32,6] = 0
$data[32,7] = 72
$data[33,0] = 'Luxemburgo'
$data[33,1] = 2487
$data[33,2] = 0
$data[33,3] = 80
$data[33,4] = 2377
$data[33,5] = 31
$data[33,6] = 0
$data[33,7] = 30
$data[34,0] = 'Pakistan'
$data[34,1] = 2458
$data[34,2] = 37
$data[34,3] = 126
$data[34,4] = 2297
$data[34,5] = 10
$data[34,6] = 1
$data[34,7] = 35
$data[35,0] = 'Indonesia'
$data[35,1] = 1986
$data[35,2] = 196
$data[35,3] = 134
$data[35,4] = 1671
$data[35,5] = 0
$data[35,6] = 11
$data[35,7] = 181
$data[36,0] = 'Tailandia'
$data[36,1] = 1978
$data[36,2] = 103
$data[36,3] = 581
$data[36,4] = 1378
$data[36,5] = 23
$data[36,6] = 4
$data[36,7] = 19
$data[37,0] = 'Arabia Saudita'
$data[37,1] = 1885
$data[37,2] = 0
$data[37,3] = 328
$data[37,4] = 1536
$data[37,5] = 31
$data[37,6] = 0
$data[37,7] = 21
$data[38,0] = 'Finlandia'
$data[38,1] = 1615
$data[38,2] = 97
$data[38,3] = 300
$data[38,4] = 1296
$data[38,5] = 65
$data[38,6] = 0
$data[38,7] = 19
$data[39,0] = 'Grecia'
$data[39,1] = 1544
$data[39,2] = 0
$data[39,3] = 61
$data[39,4] = 1430
$data[39,5] = 91
$data[39,6] = 0
$data[39,7] = 53
$data[40,0] = 'Mexico'
$data[40,1] = 1510
$data[40,2] = 132
$data[40,3] = 633
$data[40,4] = 827
$data[40,5] = 1
$data[40,6] = 13
$data[40,7] = 50
$data[41,0] = 'Panama'
$data[41,1] = 1475
$data[41,2] = 0
$data[41,3] = 9
$data[41,4] = 1429
$data[41,5] = 50
$data[41,6] = 0
$data[41,7] = 37
$data[42,0] = 'Sudafrica'
$data[42,1] = 1462
$data[42,2] = 0
$data[42,3] = 95
$data[42,4] = 1362
$data[42,5] = 7
$data[42,6] = 0
$data[42,7] = 5
$data[43,0] = 'Peru'
$data[43,1] = 1414
$data[43,2] = 0
$data[43,3] = 537
$data[43,4] = 822
$data[43,5] = 51
$data[43,6] = 0
$data[43,7] = 55
$data[44,0] = 'Republica Dominicana'
$data[44,1] = 1380
$data[44,2] = 0
$data[44,3] = 16
$data[44,4] = 1304
$data[44,5] = 147
$data[44,6] = 0
$data[44,7] = 60
$data[45,0] = 'Islandia'
$data[45,1] = 1319
$data[45,2] = 0
$data[45,3] = 284
$data[45,4] = 1031
$data[45,5] = 12
$data[45,6] = 0
$data[45,7] = 4
$data[46,0] = 'Argentina'
$data[46,1] = 1265
$data[46,2] = 0
$data[46,3] = 256
$data[46,4] = 972
$data[46,5] = 0
$data[46,6] = 1
$data[46,7] = 37
$data[47,0] = 'Serbia'
$data[47,1] = 1171
$data[47,2] = 0
$data[47,3] = 42
$data[47,4] = 1098
$data[47,5] = 81
$data[47,6] = 0
$data[47,7] = 31
$data[48,0] = 'Colombia'
$data[48,1] = 1161
$data[48,2] = 0
$data[48,3] = 55
$data[48,4] = 1087
$data[48,5] = 50
$data[48,6] = 0
$data[48,7] = 19
$data[49,0] = 'Singapur'
$data[49,1] = 1114
$data[49,2] = 65
$data[49,3] = 266
$data[49,4] = 843
$data[49,5] = 24
$data[49,6] = 1
$data[49,7] = 5
$data[50,0] = 'Emiratos Arabes Unidos'
$data[50,1] = 1024
$data[50,2] = 0
$data[50,3] = 96
$data[50,4] = 920
$data[50,5] = 2
$data[50,6] = 0
$data[50,7] = 8
$data[51,0] = 'Croacia'
$data[51,1] = 1011
$data[51,2] = 0
$data[51,3] = 88
$data[51,4] = 916
$data[51,5] = 34
$data[51,6] = 0
$data[51,7] = 7
$data[52,0] = 'Argelia'
$data[52,1] = 986
$data[52,2] = 0
$data[52,3] = 61
$data[52,4] = 839
$data[52,5] = 0
$data[52,6] = 0
$data[52,7] = 86
$data[53,0] = 'Estonia'
$data[53,1] = 961
$data[53,2] = 103
$data[53,3] = 48
$data[53,4] = 901
$data[53,5] = 16
$data[53,6] = 1
$data[53,7] = 12
$data[54,0] = 'Catar'
$data[54,1] = 949
$data[54,2] = 0
$data[54,3] = 72
$data[54,4] = 874
$data[54,5] = 37
$data[54,6] = 0
$data[54,7] = 3
$data[55,0] = 'Ucrania'
$data[55,1] = 942
$data[55,2] = 45
$data[55,3] = 19
$data[55,4] = 900
$data[55,5] = 16
$data[55,6] = 1
$data[55,7] = 23
$data[56,0] = 'Eslovenia'
$data[56,1] = 934
$data[56,2] = 37
$data[56,3] = 70
$data[56,4] = 844
$data[56,5] = 31
$data[56,6] = 3
$data[56,7] = 20
$data[57,0] = 'Nueva Zelanda'
$data[57,1] = 868
$data[57,2] = 71
$data[57,3] = 103
$data[57,4] = 764
$data[57,5] = 2
$data[57,6] = 0
$data[57,7] = 1
$data[58,0] = 'Egipto'
$data[58,1] = 865
$data[58,2] = 0
$data[58,3] = 201
$data[58,4] = 606
$data[58,5] = 0
$data[58,6] = 0
$data[58,7] = 58
$data[59,0] = 'Hong Kong'
$data[59,1] = 845
$data[59,2] = 43
$data[59,3] = 173
$data[59,4] = 668
$data[59,5] = 8
$data[59,6] = 0
$data[59,7] = 4
$data[60,0] = 'Irak'
$data[60,1] = 772
$data[60,2] = 0
$data[60,3] = 202
$data[60,4] = 516
$data[60,5] = 0
$data[60,6] = 0
$data[60,7] = 54
$data[61,0] = 'Armenia'
$data[61,1] = 736
$data[61,2] = 73
$data[61,3] = 43
$data[61,4] = 686
$data[61,5] = 30
$data[61,6] = 0
$data[61,7] = 7
$data[62,0] = 'Marruecos'
$data[62,1] = 735
$data[62,2] = 27
$data[62,3] = 49
$data[62,4] = 639
$data[62,5] = 1
$data[62,6] = 3
$data[62,7] = 47
$data[63,0] = 'Crucero'
$data[63,1] = 712
$data[63,2] = 0
$data[63,3] = 619
$data[63,4] = 82
$data[63,5] = 10
$data[63,6] = 0
$data[63,7] = 11
$data[64,0] = 'Lituania'
$data[64,1] = 696
$data[64,2] = 47
$data[64,3] = 7
$data[64,4] = 680
$data[64,5] = 11
$data[64,6] = 0
$data[64,7] = 9
$data[65,0] = 'Barein'
$data[65,1] = 672
$data[65,2] = 29
$data[65,3] = 382
$data[65,4] = 286
$data[65,5] = 3
$data[65,6] = 0
$data[65,7] = 4
$data[66,0] = 'Hungria'
$data[66,1] = 623
$data[66,2] = 38
$data[66,3] = 43
$data[66,4] = 554
$data[66,5] = 17
$data[66,6] = 5
$data[66,7] = 26
$data[67,0] = 'Bosnia y Herzegovina'
$data[67,1] = 543
$data[67,2] = 10
$data[67,3] = 27
$data[67,4] = 500
$data[67,5] = 4
$data[67,6] = 0
$data[67,7] = 16
$data[68,0] = 'Libano'
$data[68,1] = 508
$data[68,2] = 14
$data[68,3] = 46
$data[68,4] = 445
$data[68,5] = 2
$data[68,6] = 1
$data[68,7] = 17
$data[69,0] = 'Moldavia'
$data[69,1] = 505
$data[69,2] = 0
$data[69,3] = 23
$data[69,4] = 476
$data[69,5] = 65
$data[69,6] = 0
$data[69,7] = 6
$data[70,0] = 'Letonia'
$data[70,1] = 493
$data[70,2] = 35
$data[70,3] = 1
$data[70,4] = 491
$data[70,5] = 3
$data[70,6] = 1
$data[70,7] = 1
$data[71,0] = 'Bulgaria'
$data[71,1] = 477
$data[71,2] = 20
$data[71,3] = 30
$data[71,4] = 435
$data[71,5] = 18
$data[71,6] = 2
$data[71,7] = 12
$data[72,0] = 'Tunez'
$data[72,1] = 455
$data[72,2] = 0
$data[72,3] = 5
$data[72,4] = 436
$data[72,5] = 10
$data[72,6] = 0
$data[72,7] = 14
$data[73,0] = 'Kazajistan'
$data[73,1] = 453
$data[73,2] = 18
$data[73,3] = 29
$data[73,4] = 421
$data[73,5] = 6
$data[73,6] = 0
$data[73,7] = 3
$data[74,0] = 'Eslovaquia'
$data[74,1] = 450
$data[74,2] = 24
$data[74,3] = 10
$data[74,4] = 439
$data[74,5] = 3
$data[74,6] = 0
$data[74,7] = 1
$data[75,0] = 'Azerbaiyan'
$data[75,1] = 443
$data[75,2] = 43
$data[75,3] = 32
$data[75,4] = 406
$data[75,5] = 7
$data[75,6] = 0
$data[75,7] = 5
$data[76,0] = 'Principado de Andorra'
$data[76,1] = 428
$data[76,2] = 0
$data[76,3] = 10
$data[76,4] = 403
$data[76,5] = 12
$data[76,6] = 0
$data[76,7] = 15
$data[77,0] = 'Kuwait'
$data[77,1] = 417
$data[77,2] = 75
$data[77,3] = 82
$data[77,4] = 335
$data[77,5] = 16
$data[77,6] = 0
$data[77,7] = 0
$data[78,0] = 'Costa Rica'
$data[78,1] = 396
$data[78,2] = 0
$data[78,3] = 6
$data[78,4] = 388
$data[78,5] = 11
$data[78,6] = 0
$data[78,7] = 2
$data[79,0] = 'Republica de Macedonia'
$data[79,1] = 384
$data[79,2] = 0
$data[79,3] = 17
$data[79,4] = 356
$data[79,5] = 8
$data[79,6] = 0
$data[79,7] = 11
$data[80,0] = 'Uruguay'
$data[80,1] = 369
$data[80,2] = 19
$data[80,3] = 68
$data[80,4] = 297
$data[80,5] = 13
$data[80,6] = 0
$data[80,7] = 4
$data[81,0] = 'Republica de Chipre'
$data[81,1] = 356
$data[81,2] = 0
$data[81,3] = 28
$data[81,4] = 318
$data[81,5] = 11
$data[81,6] = 0
$data[81,7] = 10
$data[82,0] = 'Bielorrusia'
$data[82,1] = 351
$data[82,2] = 47
$data[82,3] = 53
$data[82,4] = 294
$data[82,5] = 11
$data[82,6] = 0
$data[82,7] = 4
$data[83,0] = 'Taiwan'
$data[83,1] = 348
$data[83,2] = 9
$data[83,3] = 50
$data[83,4] = 293
$data[83,5] = 0
$data[83,6] = 0
$data[83,7] = 5
$data[84,0] = 'Reunion'
$data[84,1] = 308
$data[84,2] = 0
$data[84,3] = 40
$data[84,4] = 268
$data[84,5] = 3
$data[84,6] = 0
$data[84,7] = 0
$data[85,0] = 'Camerun'
$data[85,1] = 306
$data[85,2] = 0
$data[85,3] = 10
$data[85,4] = 289
$data[85,5] = 0
$data[85,6] = 0
$data[85,7] = 7
$data[86,0] = 'Albania'
$data[86,1] = 304
$data[86,2] = 27
$data[86,3] = 89
$data[86,4] = 199
$data[86,5] = 7
$data[86,6] = 0
$data[86,7] = 16
$data[87,0] = 'Jordania'
$data[87,1] = 299
$data[87,2] = 0
$data[87,3] = 45
$data[87,4] = 249
$data[87,5] = 5
$data[87,6] = 0
$data[87,7] = 5
$data[88,0] = 'Burkina Faso'
$data[88,1] = 288
$data[88,2] = 0
$data[88,3] = 50
$data[88,4] = 222
$data[88,5] = 0
$data[88,6] = 0
$data[88,7] = 16
$data[89,0] = 'Afganistan'
$data[89,1] = 273
$data[89,2] = 0
$data[89,3] = 10
$data[89,4] = 257
$data[89,5] = 0
$data[89,6] = 0
$data[89,7] = 6
$data[90,0] = 'Oman'
$data[90,1] = 252
$data[90,2] = 21
$data[90,3] = 57
$data[90,4] = 194
$data[90,5] = 3
$data[90,6] = 0
$data[90,7] = 1
$data[91,0] = 'San Marino'
$data[91,1] = 245
$data[91,2] = 0
$data[91,3] = 21
$data[91,4] = 194
$data[91,5] = 15
$data[91,6] = 0
$data[91,7] = 30
$data[92,0] = 'Vietnam'
$data[92,1] = 237
$data[92,2] = 4
$data[92,3] = 85
$data[92,4] = 152
$data[92,5] = 3
$data[92,6] = 0
$data[92,7] = 0
$data[93,0] = 'Cuba'
$data[93,1] = 233
$data[93,2] = 0
$data[93,3] = 13
$data[93,4] = 214
$data[93,5] = 7
$data[93,6] = 0
$data[93,7] = 6
$data[94,0] = 'Honduras'
$data[94,1] = 222
$data[94,2] = 3
$data[94,3] = 3
$data[94,4] = 204
$data[94,5] = 10
$data[94,6] = 1
$data[94,7] = 15
$data[95,0] = 'Uzbekistan'
$data[95,1] = 221
$data[95,2] = 16
$data[95,3] = 25
$data[95,4] = 194
$data[95,5] = 8
$data[95,6] = 0
$data[95,7] = 2
$data[96,0] = 'Senegal'
$data[96,1] = 207
$data[96,2] = 12
$data[96,3] = 66
$data[96,4] = 140
$data[96,5] = 1
$data[96,6] = 0
$data[96,7] = 1
$data[97,0] = 'Ghana'
$data[97,1] = 204
$data[97,2] = 0
$data[97,3] = 31
$data[97,4] = 168
$data[97,5] = 2
$data[97,6] = 0
$data[97,7] = 5
$data[98,0] = 'Malta'
$data[98,1] = 202
$data[98,2] = 6
$data[98,3] = 2
$data[98,4] = 200
$data[98,5] = 2
$data[98,6] = 0
$data[98,7] = 0
$data[99,0] = 'Costa de Marfil'
$data[99,1] = 194
$data[99,2] = 0
$data[99,3] = 15
$data[99,4] = 178
$data[99,5] = 0
$data[99,6] = 0
$data[99,7] = 1
$data[100,0] = 'Nigeria'
$data[100,1] = 190
$data[100,2] = 6
$data[100,3] = 20
$data[100,4] = 168
$data[100,5] = 0
$data[100,6] = 0
$data[100,7] = 2
$data[101,0] = 'Islas Feroe'
$data[101,1] = 179
$data[101,2] = 2
$data[101,3] = 91
$data[101,4] = 88
$data[101,5] = 1
$data[101,6] = 0
$data[101,7] = 0
$data[102,0] = 'Estado de Palestina'
$data[102,1] = 171
$data[102,2] = 10
$data[102,3] = 18
$data[102,4] = 152
$data[102,5] = 0
$data[102,6] = 0
$data[102,7] = 1
$data[103,0] = 'Mauricio'
$data[103,1] = 169
$data[103,2] = 0
$data[103,3] = 0
$data[103,4] = 162
$data[103,5] = 1
$data[103,6] = 0
$data[103,7] = 7
$data[104,0] = 'Montenegro'
$data[104,1] = 160
$data[104,2] = 16
$data[104,3] = 0
$data[104,4] = 158
$data[104,5] = 4
$data[104,6] = 0
$data[104,7] = 2
$data[105,0] = 'Sri Lanka'
$data[105,1] = 151
$data[105,2] = 0
$data[105,3] = 22
$data[105,4] = 125
$data[105,5] = 5
$data[105,6] = 0
$data[105,7] = 4
$data[106,0] = 'Georgia'
$data[106,1] = 148
$data[106,2] = 14
$data[106,3] = 27
$data[106,4] = 121
$data[106,5] = 6
$data[106,6] = 0
$data[106,7] = 0
$data[107,0] = 'Venezuela'
$data[107,1] = 146
$data[107,2] = 0
$data[107,3] = 43
$data[107,4] = 98
$data[107,5] = 6
$data[107,6] = 0
$data[107,7] = 5
$data[108,0] = 'Martinica'
$data[108,1] = 138
$data[108,2] = 0
$data[108,3] = 27
$data[108,4] = 108
$data[108,5] = 19
$data[108,6] = 0
$data[108,7] = 3
$data[109,0] = 'Consejo Danes para los Refugiados'
$data[109,1] = 134
$data[109,2] = 0
$data[109,3] = 3
$data[109,4] = 118
$data[109,5] = 0
$data[109,6] = 0
$data[109,7] = 13
$data[110,0] = 'Brunei'
$data[110,1] = 134
$data[110,2] = 1
$data[110,3] = 65
$data[110,4] = 68
$data[110,5] = 3
$data[110,6] = 0
$data[110,7] = 1
$data[111,0] = 'Bolivia'
$data[111,1] = 132
$data[111,2] = 9
$data[111,3] = 1
$data[111,4] = 122
$data[111,5] = 3
$data[111,6] = 1
$data[111,7] = 9
$data[112,0] = 'Kirguistan'
$data[112,1] = 130
$data[112,2] = 14
$data[112,3] = 5
$data[112,4] = 124
$data[112,5] = 5
$data[112,6] = 1
$data[112,7] = 1
$data[113,0] = 'Guadalupe'
$data[113,1] = 128
$data[113,2] = 0
$data[113,3] = 24
$data[113,4] = 98
$data[113,5] = 14
$data[113,6] = 0
$data[113,7] = 6
$data[114,0] = 'Mayotte'
$data[114,1] = 116
$data[114,2] = 0
$data[114,3] = 10
$data[114,4] = 105
$data[114,5] = 3
$data[114,6] = 0
$data[114,7] = 1
$data[115,0] = 'Isla de Man'
$data[115,1] = 114
$data[115,2] = 19
$data[115,3] = 0
$data[115,4] = 113
$data[115,5] = 0
$data[115,6] = 0
$data[115,7] = 1
$data[116,0] = 'Camboya'
$data[116,1] = 114
$data[116,2] = 4
$data[116,3] = 35
$data[116,4] = 79
$data[116,5] = 1
$data[116,6] = 0
$data[116,7] = 0
$data[117,0] = 'Kenia'
$data[117,1] = 110
$data[117,2] = 0
$data[117,3] = 4
$data[117,4] = 103
$data[117,5] = 2
$data[117,6] = 0
$data[117,7] = 3
$data[118,0] = 'Niger'
$data[118,1] = 98
$data[118,2] = 0
$data[118,3] = 0
$data[118,4] = 93
$data[118,5] = 0
$data[118,6] = 0
$data[118,7] = 5
$data[119,0] = 'Trinidad yTobago'
$data[119,1] = 97
$data[119,2] = 3
$data[119,3] = 1
$data[119,4] = 90
$data[119,5] = 0
$data[119,6] = 1
$data[119,7] = 6
$data[120,0] = 'Paraguay'
$data[120,1] = 92
$data[120,2] = 15
$data[120,3] = 4
$data[120,4] = 85
$data[120,5] = 4
$data[120,6] = 0
$data[120,7] = 3
$data[121,0] = 'Gibraltar'
$data[121,1] = 88
$data[121,2] = 0
$data[121,3] = 46
$data[121,4] = 42
$data[121,5] = 0
$data[121,6] = 0
$data[121,7] = 0
$data[122,0] = 'Ruanda'
$data[122,1] = 84
$data[122,2] = 0
$data[122,3] = 0
$data[122,4] = 84
$data[122,5] = 0
$data[122,6] = 0
$data[122,7] = 0
$data[123,0] = 'Liechtenstein'
$data[123,1] = 75
$data[123,2] = 0
$data[123,3] = 0
$data[123,4] = 75
$data[123,5] = 0
$data[123,6] = 0
$data[123,7] = 0
$data[124,0] = 'Banglades'
$data[124,1] = 61
$data[124,2] = 5
$data[124,3] = 26
$data[124,4] = 29
$data[124,5] = 1
$data[124,6] = 0
$data[124,7] = 6
$data[125,0] = 'Aruba'
$data[125,1] = 60
$data[125,2] = 0
$data[125,3] = 1
$data[125,4] = 59
$data[125,5] = 0
$data[125,6] = 0
$data[125,7] = 0
$data[126,0] = 'Monaco'
$data[126,1] = 60
$data[126,2] = 0
$data[126,3] = 2
$data[126,4] = 57
$data[126,5] = 2
$data[126,6] = 0
$data[126,7] = 1
$data[127,0] = 'Madagascar'
$data[127,1] = 59
$data[127,2] = 0
$data[127,3] = 0
$data[127,4] = 59
$data[127,5] = 6
$data[127,6] = 0
$data[127,7] = 0
$data[128,0] = 'Guinea'
$data[128,1] = 52
$data[128,2] = 0
$data[128,3] = 0
$data[128,4] = 52
$data[128,5] = 0
$data[128,6] = 0
$data[128,7] = 0
$data[129,0] = 'Guayana Francesa'
$data[129,1] = 51
$data[129,2] = 0
$data[129,3] = 15
$data[129,4] = 36
$data[129,5] = 0
$data[129,6] = 0
$data[129,7] = 0
$data[130,0] = 'Republica de Yibuti'
$data[130,1] = 49
$data[130,2] = 9
$data[130,3] = 8
$data[130,4] = 41
$data[130,5] = 0
$data[130,6] = 0
$data[130,7] = 0
$data[131,0] = 'Guatemala'
$data[131,1] = 47
$data[131,2] = 0
$data[131,3] = 2
$data[131,4] = 42
$data[131,5] = 0
$data[131,6] = 0
$data[131,7] = 3
$data[132,0] = 'Barbados'
$data[132,1] = 47
$data[132,2] = 0
$data[132,3] = 12
$data[132,4] = 34
$data[132,5] = 1
$data[132,6] = 0
$data[132,7] = 1
$data[133,0] = 'El Salvador'
$data[133,1] = 46
$data[133,2] = 0
$data[133,3] = 0
$data[133,4] = 46
$data[133,5] = 0
$data[133,6] = 0
$data[133,7] = 0
$data[134,0] = 'Uganda'
$data[134,1] = 46
$data[134,2] = 5
$data[134,3] = 0
$data[134,4] = 44
$data[134,5] = 4
$data[134,6] = 0
$data[134,7] = 2
$data[135,0] = 'Jamaica'
$data[135,1] = 45
$data[135,2] = 0
$data[135,3] = 0
$data[135,4] = 45
$data[135,5] = 0
$data[135,6] = 0
$data[135,7] = 0
$data[136,0] = 'Macao'
$data[136,1] = 42
$data[136,2] = 1
$data[136,3] = 10
$data[136,4] = 32
$data[136,5] = 0
$data[136,6] = 0
$data[136,7] = 0
$data[137,0] = 'Zambia'
$data[137,1] = 39
$data[137,2] = 0
$data[137,3] = 0
$data[137,4] = 38
$data[137,5] = 0
$data[137,6] = 0
$data[137,7] = 1
$data[138,0] = 'Puerto Rico'
$data[138,1] = 39
$data[138,2] = 0
$data[138,3] = 1
$data[138,4] = 36
$data[138,5] = 0
$data[138,6] = 0
$data[138,7] = 2
$data[139,0] = 'Togo'
$data[139,1] = 39
$data[139,2] = 0
$data[139,3] = 17
$data[139,4] = 20
$data[139,5] = 0
$data[139,6] = 0
$data[139,7] = 2
$data[140,0] = 'Polinesia Francesa'
$data[140,1] = 37
$data[140,2] = 0
$data[140,3] = 0
$data[140,4] = 37
$data[140,5] = 1
$data[140,6] = 0
$data[140,7] = 0
$data[141,0] = 'Mali'
$data[141,1] = 36
$data[141,2] = 0
$data[141,3] = 0
$data[141,4] = 33
$data[141,5] = 0
$data[141,6] = 0
$data[141,7] = 3
$data[142,0] = 'Etiopia'
$data[142,1] = 35
$data[142,2] = 6
$data[142,3] = 3
$data[142,4] = 32
$data[142,5] = 2
$data[142,6] = 0
$data[142,7] = 0
$data[143,0] = 'Bermudas'
$data[143,1] = 35
$data[143,2] = 0
$data[143,3] = 11
$data[143,4] = 24
$data[143,5] = 0
$data[143,6] = 0
$data[143,7] = 0
$data[144,0] = 'Guam'
$data[144,1] = 32
$data[144,2] = 0
$data[144,3] = 0
$data[144,4] = 31
$data[144,5] = 0
$data[144,6] = 0
$data[144,7] = 1
$data[145,0] = 'Islas Caimanes'
$data[145,1] = 28
$data[145,2] = 0
$data[145,3] = 0
$data[145,4] = 27
$data[145,5] = 0
$data[145,6] = 0
$data[145,7] = 1
$data[146,0] = 'Bahamas'
$data[146,1] = 24
$data[146,2] = 0
$data[146,3] = 1
$data[146,4] = 22
$data[146,5] = 1
$data[146,6] = 0
$data[146,7] = 1
$data[147,0] = 'San Martin (Parte Holandesa)'
$data[147,1] = 23
$data[147,2] = 5
$data[147,3] = 6
$data[147,4] = 15
$data[147,5] = 0
$data[147,6] = 1
$data[147,7] = 2
$data[148,0] = 'Eritrea'
$data[148,1] = 22
$data[148,2] = 0
$data[148,3] = 0
$data[148,4] = 22
$data[148,5] = 0
$data[148,6] = 0
$data[148,7] = 0
$data[149,0] = 'San Martin (Parte Francesa)'
$data[149,1] = 22
$data[149,2] = 0
$data[149,3] = 2
$data[149,4] = 19
$data[149,5] = 0
$data[149,6] = 0
$data[149,7] = 1
$data[150,0] = 'Congo'
$data[150,1] = 22
$data[150,2] = 0
$data[150,3] = 2
$data[150,4] = 18
$data[150,5] = 0
$data[150,6] = 0
$data[150,7] = 2
$data[151,0] = 'Gabon'
$data[151,1] = 21
$data[151,2] = 0
$data[151,3] = 1
$data[151,4] = 19
$data[151,5] = 0
$data[151,6] = 0
$data[151,7] = 1
$data[152,0] = 'Birmania'
$data[152,1] = 20
$data[152,2] = 0
$data[152,3] = 0
$data[152,4] = 19
$data[152,5] = 0
$data[152,6] = 0
$data[152,7] = 1
$data[153,0] = 'Tanzania'
$data[153,1] = 20
$data[153,2] = 0
$data[153,3] = 3
$data[153,4] = 16
$data[153,5] = 0
$data[153,6] = 0
$data[153,7] = 1
$data[154,0] = 'Guyana'
$data[154,1] = 19
$data[154,2] = 0
$data[154,3] = 0
$data[154,4] = 15
$data[154,5] = 0
$data[154,6] = 0
$data[154,7] = 4
$data[155,0] = 'Maldivas'
$data[155,1] = 19
$data[155,2] = 0
$data[155,3] = 13
$data[155,4] = 6
$data[155,5] = 0
$data[155,6] = 0
$data[155,7] = 0
$data[156,0] = 'Nueva Caledonia'
$data[156,1] = 18
$data[156,2] = 0
$data[156,3] = 1
$data[156,4] = 17
$data[156,5] = 0
$data[156,6] = 0
$data[156,7] = 0
$data[157,0] = 'Haiti'
$data[157,1] = 18
$data[157,2] = 2
$data[157,3] = 1
$data[157,4] = 17
$data[157,5] = 0
$data[157,6] = 0
$data[157,7] = 0
$data[158,0] = 'Islas Virgenes de los Estados Unidos'
$data[158,1] = 17
$data[158,2] = 0
$data[158,3] = 0
$data[158,4] = 17
$data[158,5] = 0
$data[158,6] = 0
$data[158,7] = 0
$data[159,0] = 'Siria'
$data[159,1] = 16
$data[159,2] = 0
$data[159,3] = 0
$data[159,4] = 14
$data[159,5] = 0
$data[159,6] = 0
$data[159,7] = 2
$data[160,0] = 'Guinea Ecuatorial'
$data[160,1] = 15
$data[160,2] = 0
$data[160,3] = 1
$data[160,4] = 14
$data[160,5] = 0
$data[160,6] = 0
$data[160,7] = 0
$data[161,0] = 'Mongolia'
$data[161,1] = 14
$data[161,2] = 0
$data[161,3] = 2
$data[161,4] = 12
$data[161,5] = 0
$data[161,6] = 0
$data[161,7] = 0
$data[162,0] = 'Namibia'
$data[162,1] = 14
$data[162,2] = 0
$data[162,3] = 3
$data[162,4] = 11
$data[162,5] = 0
$data[162,6] = 0
$data[162,7] = 0
$data[163,0] = 'Santa Lucia'
$data[163,1] = 13
$data[163,2] = 0
$data[163,3] = 1
$data[163,4] = 12
$data[163,5] = 0
$data[163,6] = 0
$data[163,7] = 0
$data[164,0] = 'Benin'
$data[164,1] = 13
$data[164,2] = 0
$data[164,3] = 1
$data[164,4] = 12
$data[164,5] = 0
$data[164,6] = 0
$data[164,7] = 0
$data[165,0] = 'Dominica'
$data[165,1] = 12
$data[165,2] = 0
$data[165,3] = 0
$data[165,4] = 12
$data[165,5] = 0
$data[165,6] = 0
$data[165,7] = 0
$data[166,0] = 'Libia'
$data[166,1] = 11
$data[166,2] = 0
$data[166,3] = 0
$data[166,4] = 10
$data[166,5] = 0
$data[166,6] = 0
$data[166,7] = 1
$data[167,0] = 'Curazao'
$data[167,1] = 11
$data[167,2] = 0
$data[167,3] = 3
$data[167,4] = 7
$data[167,5] = 0
$data[167,6] = 0
$data[167,7] = 1
$data[168,0] = 'Granada'
$data[168,1] = 10
$data[168,2] = 0
$data[168,3] = 0
$data[168,4] = 10
$data[168,5] = 0
$data[168,6] = 0
$data[168,7] = 0
$data[169,0] = 'Laos'
$data[169,1] = 10
$data[169,2] = 0
$data[169,3] = 0
$data[169,4] = 10
$data[169,5] = 0
$data[169,6] = 0
$data[169,7] = 0
$data[170,0] = 'Seychelles'
$data[170,1] = 10
$data[170,2] = 0
$data[170,3] = 0
$data[170,4] = 10
$data[170,5] = 0
$data[170,6] = 0
$data[170,7] = 0
$data[171,0] = 'Surinam'
$data[171,1] = 10
$data[171,2] = 0
$data[171,3] = 0
$data[171,4] = 10
$data[171,5] = 0
$data[171,6] = 0
$data[171,7] = 0
$data[172,0] = 'Mozambique'
$data[172,1] = 10
$data[172,2] = 0
$data[172,3] = 0
$data[172,4] = 10
$data[172,5] = 0
$data[172,6] = 0
$data[172,7] = 0
$data[173,0] = 'Groenlandia'
$data[173,1] = 10
$data[173,2] = 0
$data[173,3] = 3
$data[173,4] = 7
$data[173,5] = 0
$data[173,6] = 0
$data[173,7] = 0
$data[174,0] = 'Suazilandia'
$data[174,1] = 9
$data[174,2] = 0
$data[174,3] = 0
$data[174,4] = 9
$data[174,5] = 0
$data[174,6] = 0
$data[174,7] = 0
$data[175,0] = 'San Cristobal y Nieves'
$data[175,1] = 9
$data[175,2] = 0
$data[175,3] = 0
$data[175,4] = 9
$data[175,5] = 0
$data[175,6] = 0
$data[175,7] = 0
$data[176,0] = 'Guinea-Bisau'
$data[176,1] = 9
$data[176,2] = 0
$data[176,3] = 0
$data[176,4] = 9
$data[176,5] = 0
$data[176,6] = 0
$data[176,7] = 0
$data[177,0] = 'Antigua y Barbuda'
$data[177,1] = 9
$data[177,2] = 0
$data[177,3] = 0
$data[177,4] = 9
$data[177,5] = 0
$data[177,6] = 0
$data[177,7] = 0
$data[178,0] = 'Zimbabue'
$data[178,1] = 9
$data[178,2] = 0
$data[178,3] = 0
$data[178,4] = 8
$data[178,5] = 0
$data[178,6] = 0
$data[178,7] = 1
$data[179,0] = 'Montserrat'
$data[179,1] = 9
$data[179,2] = 0
$data[179,3] = 0
$data[179,4] = 7
$data[179,5] = 0
$data[179,6] = 0
$data[179,7] = 2
$data[180,0] = 'Republica de Africa Central'
$data[180,1] = 8
$data[180,2] = 5
$data[180,3] = 0
$data[180,4] = 8
$data[180,5] = 0
$data[180,6] = 0
$data[180,7] = 0
$data[181,0] = 'Republica del Chad'
$data[181,1] = 8
$data[181,2] = 0
$data[181,3] = 0
$data[181,4] = 8
$data[181,5] = 0
$data[181,6] = 0
$data[181,7] = 0
$data[182,0] = 'Angola'
$data[182,1] = 8
$data[182,2] = 0
$data[182,3] = 1
$data[182,4] = 5
$data[182,5] = 0
$data[182,6] = 0
$data[182,7] = 2
$data[183,0] = 'Sudan'
$data[183,1] = 8
$data[183,2] = 0
$data[183,3] = 2
$data[183,4] = 4
$data[183,5] = 0
$data[183,6] = 0
$data[183,7] = 2
$data[184,0] = 'Santa Sede'
$data[184,1] = 7
$data[184,2] = 0
$data[184,3] = 0
$data[184,4] = 7
$data[184,5] = 0
$data[184,6] = 0
$data[184,7] = 0
$data[185,0] = 'Fiyi'
$data[185,1] = 7
$data[185,2] = 0
$data[185,3] = 0
$data[185,4] = 7
$data[185,5] = 0
$data[185,6] = 0
$data[185,7] = 0
$data[186,0] = 'Liberia'
$data[186,1] = 6
$data[186,2] = 0
$data[186,3] = 0
$data[186,4] = 6
$data[186,5] = 0
$data[186,6] = 0
$data[186,7] = 0
$data[187,0] = 'San Bartolome'
$data[187,1] = 6
$data[187,2] = 0
$data[187,3] = 1
$data[187,4] = 5
$data[187,5] = 0
$data[187,6] = 0
$data[187,7] = 0
$data[188,0] = 'Nepal'
$data[188,1] = 6
$data[188,2] = 0
$data[188,3] = 1
$data[188,4] = 5
$data[188,5] = 0
$data[188,6] = 0
$data[188,7] = 0
$data[189,0] = 'Cabo Verde'
$data[189,1] = 6
$data[189,2] = 0
$data[189,3] = 0
$data[189,4] = 5
$data[189,5] = 0
$data[189,6] = 0
$data[189,7] = 1
$data[190,0] = 'Mauritania'
$data[190,1] = 6
$data[190,2] = 0
$data[190,3] = 2
$data[190,4] = 3
$data[190,5] = 0
$data[190,6] = 0
$data[190,7] = 1
$data[191,0] = 'Islas Turcas y Caicos'
$data[191,1] = 5
$data[191,2] = 0
$data[191,3] = 0
$data[191,4] = 5
$data[191,5] = 0
$data[191,6] = 0
$data[191,7] = 0
$data[192,0] = 'Somalia'
$data[192,1] = 5
$data[192,2] = 0
$data[192,3] = 1
$data[192,4] = 4
$data[192,5] = 0
$data[192,6] = 0
$data[192,7] = 0
$data[193,0] = 'Nicaragua'
$data[193,1] = 5
$data[193,2] = 0
$data[193,3] = 0
$data[193,4] = 4
$data[193,5] = 0
$data[193,6] = 0
$data[193,7] = 1
$data[194,0] = 'Butan'
$data[194,1] = 5
$data[194,2] = 0
$data[194,3] = 2
$data[194,4] = 3
$data[194,5] = 0
$data[194,6] = 0
$data[194,7] = 0
$data[195,0] = 'Botsuana'
$data[195,1] = 4
$data[195,2] = 0
$data[195,3] = 0
$data[195,4] = 3
$data[195,5] = 0
$data[195,6] = 0
$data[195,7] = 1
$data[196,0] = 'Gambia'
$data[196,1] = 4
$data[196,2] = 0
$data[196,3] = 2
$data[196,4] = 1
$data[196,5] = 0
$data[196,6] = 0
$data[196,7] = 1
$data[197,0] = 'Belice'
$data[197,1] = 3
$data[197,2] = 0
$data[197,3] = 0
$data[197,4] = 3
$data[197,5] = 0
$data[197,6] = 0
$data[197,7] = 0
$data[198,0] = 'Anguila'
$data[198,1] = 3
$data[198,2] = 0
$data[198,3] = 0
$data[198,4] = 3
$data[198,5] = 0
$data[198,6] = 0
$data[198,7] = 0
$data[199,0] = 'Islas Virgenes Britanicas'
$data[199,1] = 3
$data[199,2] = 0
$data[199,3] = 0
$data[199,4] = 3
$data[199,5] = 0
$data[199,6] = 0
$data[199,7] = 0
$data[200,0] = 'Malaui'
$data[200,1] = 3
$data[200,2] = 0
$data[200,3] = 0
$data[200,4] = 3
$data[200,5] = 0
$data[200,6] = 0
$data[200,7] = 0
$data[201,0] = 'Burundi'
$data[201,1] = 3
$data[201,2] = 0
$data[201,3] = 0
$data[201,4] = 3
$data[201,5] = 0
$data[201,6] = 0
$data[201,7] = 0
$data[202,0] = 'Sierra Leona'
$data[202,1] = 2
$data[202,2] = 0
$data[202,3] = 0
$data[202,4] = 2
$data[202,5] = 0
$data[202,6] = 0
$data[202,7] = 0
$data[203,0] = 'Bonaire, San Eustaquio y Saba'
$data[203,1] = 2
$data[203,2] = 0
$data[203,3] = 0
$data[203,4] = 2
$data[203,5] = 0
$data[203,6] = 0
$data[203,7] = 0
$data[204,0] = 'San Vicente y las Granadinas'
$data[204,1] = 2
$data[204,2] = 0
$data[204,3] = 1
$data[204,4] = 1
$data[204,5] = 0
$data[204,6] = 0
$data[204,7] = 0
$data[205,0] = 'Timor Oriental'
$data[205,1] = 1
$data[205,2] = 0
$data[205,3] = 0
$data[205,4] = 1
$data[205,5] = 0
$data[205,6] = 0
$data[205,7] = 0
$data[206,0] = 'Papua Nueva Guinea'
$data[206,1] = 1
$data[206,2] = 0
$data[206,3] = 0
$data[206,4] = 1
$data[206,5] = 0
$data[206,6] = 0
$data[206,7] = 0

$ws.Range("A4:H210").Value = $data

$ws.Range("A1").Value = 'Datos actualizados a 3 de Abril de 2020 a las 13:50'
